# Auto-generated edit script: applies value updates to Ifrit_Profits-equivalent sheets
# (workbook tabs ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 190
$ws.Range("I58").Value = 190
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 570
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -420
$ws.Range("N58").ClearContents()
$ws.Range("H64").Value = 4328.421
$ws.Range("I64").Value = 4234
$ws.Range("J64").Value = 4433.3335
$ws.Range("K64").Value = 4234
$ws.Range("L64").Value = 4433.3335
$ws.Range("M64").Value = -3986
$ws.Range("N64").Value = -4929.3335
$ws.Range("H67").Value = 4328.421
$ws.Range("I67").Value = 4234
$ws.Range("J67").Value = 4433.3335
$ws.Range("K67").Value = 4234
$ws.Range("L67").Value = 4433.3335
$ws.Range("M67").Value = -3376
$ws.Range("N67").Value = -6149.3335
$ws.Range("H80").Value = 753.1875
$ws.Range("I80").Value = 742.46155
$ws.Range("J80").Value = 799.6667
$ws.Range("K80").Value = 2227.38465
$ws.Range("L80").Value = 2399.0001
$ws.Range("M80").Value = -1229.38465
$ws.Range("N80").Value = -4395.0001
$ws.Range("H82").Value = 2500
$ws.Range("I82").Value = 2500
$ws.Range("K82").Value = 7500
$ws.Range("M82").Value = -7094
$ws.Range("H83").Value = 753.1875
$ws.Range("I83").Value = 742.46155
$ws.Range("J83").Value = 799.6667
$ws.Range("K83").Value = 6682.15395
$ws.Range("L83").Value = 7197.0003
$ws.Range("M83").Value = -1690.15395
$ws.Range("N83").Value = -17181.0003
$ws.Range("H85").Value = 2500
$ws.Range("I85").Value = 2500
$ws.Range("K85").Value = 7500
$ws.Range("M85").Value = -6096
$ws.Range("H86").Value = 3668.182
$ws.Range("I86").Value = 2130.6
$ws.Range("J86").Value = 4949.5
$ws.Range("K86").Value = 2130.6
$ws.Range("L86").Value = 4949.5
$ws.Range("M86").Value = -1007.6
$ws.Range("N86").Value = -7195.5
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42496
$ws.Range("H88").Value = 1907.5
$ws.Range("I88").Value = 650
$ws.Range("J88").Value = 2326.6667
$ws.Range("K88").Value = 650
$ws.Range("L88").Value = 2326.6667
$ws.Range("M88").Value = -244
$ws.Range("N88").Value = -3138.6667
$ws.Range("H89").Value = 3668.182
$ws.Range("I89").Value = 2130.6
$ws.Range("J89").Value = 4949.5
$ws.Range("K89").Value = 10653
$ws.Range("L89").Value = 24747.5
$ws.Range("M89").Value = -5037
$ws.Range("N89").Value = -35979.5
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -132480
$ws.Range("H91").Value = 1907.5
$ws.Range("I91").Value = 650
$ws.Range("J91").Value = 2326.6667
$ws.Range("K91").Value = 650
$ws.Range("L91").Value = 2326.6667
$ws.Range("M91").Value = 754
$ws.Range("N91").Value = -5134.6667
$ws.Range("H118").Value = 423.76923
$ws.Range("I118").Value = 423.76923
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1271.30769
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 385.6923099999999
$ws.Range("N118").ClearContents()
$ws.Range("H137").Value = 52633588
$ws.Range("I137").Value = 1820.2
$ws.Range("K137").Value = 5460.6
$ws.Range("M137").Value = -2910.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12956
$ws.Range("I32").Value = 10534.052
$ws.Range("K32").Value = 10534.052
$ws.Range("M32").Value = -10247.052
$ws.Range("H45").Value = 1534.2273
$ws.Range("I45").Value = 1262.75
$ws.Range("K45").Value = 1262.75
$ws.Range("M45").Value = -885.75
$ws.Range("H61").Value = 4275639
$ws.Range("I61").Value = 4832920
$ws.Range("K61").Value = 4832920
$ws.Range("M61").Value = -4832708
$ws.Range("H136").Value = 4275639
$ws.Range("I136").Value = 4832920
$ws.Range("K136").Value = 14498760
$ws.Range("M136").Value = -14496210

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 12570002
$ws.Range("I134").Value = 13407890
$ws.Range("J134").Value = 1680
$ws.Range("K134").Value = 40223670
$ws.Range("L134").Value = 5040
$ws.Range("M134").Value = -40221135
$ws.Range("N134").Value = -10110

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 28159.29
$ws.Range("I99").Value = 2011.4584
$ws.Range("K99").Value = 2011.4584
$ws.Range("M99").Value = -513.4584
$ws.Range("N99").Value = -75980.14
$ws.Range("H126").Value = 28159.29
$ws.Range("I126").Value = 2011.4584
$ws.Range("K126").Value = 6034.3752
$ws.Range("M126").Value = -3564.3752
$ws.Range("N126").Value = -223892.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1030.1111
$ws.Range("I44").Value = 574.2
$ws.Range("J44").Value = 1600
$ws.Range("K44").Value = 1722.6
$ws.Range("L44").Value = 4800
$ws.Range("M44").Value = -1324.6
$ws.Range("N44").Value = -5596
$ws.Range("H46").Value = 2326.3157
$ws.Range("I46").Value = 1525
$ws.Range("J46").Value = 2909.0908
$ws.Range("K46").Value = 4575
$ws.Range("L46").Value = 8727.2724
$ws.Range("M46").Value = -4484
$ws.Range("N46").Value = -8909.2724
$ws.Range("H69").Value = 1198.8889
$ws.Range("I69").Value = 350
$ws.Range("J69").Value = 1441.4286
$ws.Range("K69").Value = 1050
$ws.Range("L69").Value = 4324.2858
$ws.Range("M69").Value = -239
$ws.Range("N69").Value = -5946.2858
$ws.Range("H72").Value = 1198.8889
$ws.Range("I72").Value = 350
$ws.Range("J72").Value = 1441.4286
$ws.Range("K72").Value = 3150
$ws.Range("L72").Value = 12972.8574
$ws.Range("M72").Value = 906
$ws.Range("N72").Value = -21084.8574
$ws.Range("H138").Value = 4896.6665
$ws.Range("I138").Value = 3752
$ws.Range("J138").Value = 5714.2856
$ws.Range("K138").Value = 11256
$ws.Range("L138").Value = 17142.8568
$ws.Range("M138").Value = -6116
$ws.Range("N138").Value = -27422.8568
$ws.Range("H141").Value = 5938.8887
$ws.Range("I141").Value = 5938.8887
$ws.Range("K141").Value = 17816.6661
$ws.Range("M141").Value = -12636.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 43851
$ws.Range("J110").Value = 43851
$ws.Range("L110").Value = 43851
$ws.Range("N110").Value = -52031
$ws.Range("H113").Value = 1612.0769
$ws.Range("I113").Value = 1333
$ws.Range("J113").Value = 2240
$ws.Range("K113").Value = 1333
$ws.Range("L113").Value = 2240
$ws.Range("M113").Value = 837
$ws.Range("N113").Value = -6580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 348.9524
$ws.Range("I113").Value = 351.15384
$ws.Range("J113").Value = 345.375
$ws.Range("K113").Value = 1053.46152
$ws.Range("L113").Value = 1036.125
$ws.Range("M113").Value = 1116.53848
$ws.Range("N113").Value = -5376.125
